$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "router3"
$ws.Range("B4").Value = "third_device"

$ws.Range("B1:B4").Select()
